# Bugfix: remove the first 16 stale QoQ GDP observations (rows 2-17) that
# were evaluated against an outdated release window. Deleting these rows
# shifts the remaining observations up so the series again starts at the
# correct first release date, shrinking the used range from A1:B164 to
# A1:B148.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("2:17").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp) | Out-Null
